$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.923.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7746"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3129"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07247"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08716"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.079.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7721"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.410"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.215"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.912.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.333.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007886"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.177"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1591"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.523"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.048"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.544"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.522"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.126"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05439"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.249"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7551"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.65%  "

$ws.Range("E37").Value = "  +0.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01987"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.784"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4520"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.093"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.14%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.098.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.21%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.239.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.44%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8544"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.887"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.617"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.831"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
